# ---------------------------------------------------------------------------
# Edit: rename the stat sheets to human-friendly titles, and bump every
# player's "Age" column (format YY-DDD, years-days) forward by one day to
# reflect the stats being (re)generated a day later.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------------
$renames = @{
    "StandardStats"     = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $wb.Worksheets.Item($oldName).Name = $renames[$oldName]
}

# --- 2. Bump the Age column (column E) forward by one day ------------------
# Every stats sheet (all except "Matches") carries the same per-player Age
# values in column E, formatted as "YY-DDD" (years-days). The new snapshot
# increments the day-of-year component by 1 for each populated row.

$ageUpdates = @(
    @(4,  "23-288"), @(5,  "24-037"), @(6,  "31-125"), @(7,  "25-298"),
    @(8,  "25-043"), @(9,  "32-063"), @(10, "23-248"), @(11, "26-064"),
    @(12, "33-064"), @(13, "26-132"), @(14, "31-201"), @(15, "26-267"),
    @(16, "34-319"), @(17, "28-112"), @(18, "31-361"), @(19, "27-277"),
    @(20, "21-257"), @(21, "32-314"), @(22, "29-054"), @(23, "23-323"),
    @(24, "28-272"), @(25, "23-155"), @(26, "21-002"), @(27, "18-263"),
    @(28, "29-151"), @(29, "22-126"), @(30, "25-044"), @(31, "21-020")
)

# "StandardStats" and "PlayingTime" have a few extra populated rows beyond
# row 31 that the other stat sheets don't have.
$extraAgeUpdates = @(
    @(32, "31-251"), @(33, "21-161"), @(34, "22-329"), @(35, "18-307"),
    @(37, "19-002"), @(38, "16-313"), @(39, "23-213")
)

$statSheetNames = @(
    "Standard Stats", "Shooting Stats", "Passing Stats", "Pass Types",
    "Goal & Shot Creation", "Defensive Actions", "Possession",
    "Playing Time", "Miscellaneous Stats"
)

$extendedSheetNames = @("Standard Stats", "Playing Time")

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($pair in $ageUpdates) {
        $row = $pair[0]
        $newAge = $pair[1]
        $ws.Range("E$row").Value = $newAge
    }

    if ($extendedSheetNames -contains $sheetName) {
        foreach ($pair in $extraAgeUpdates) {
            $row = $pair[0]
            $newAge = $pair[1]
            $ws.Range("E$row").Value = $newAge
        }
    }
}
